# secretaires.xlsx — "Add files via upload / new version fichier secretaires.csv"
#
# Adds a new "Pays" (country) column (F) with the value "British Isles" for
# every existing data row, and appends a new record (row 118) for the PDF
# "89.pdf" / "The Ladies and Gentlemen's Complete Letter-Writer ..." by
# Anon., dated 1797, also tagged "British Isles".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append the new row (118) -------------------------------------------
# Pick up the italic "Anon." title style (col B) and the non-italic author
# style (col C) from the row directly above, so the new record matches the
# existing "Anon." entries' look.
$ws.Cells.Item(117, 2).Copy()
$ws.Cells.Item(118, 2).PasteSpecial(-4122)   # xlPasteFormats
$ws.Cells.Item(117, 3).Copy()
$ws.Cells.Item(118, 3).PasteSpecial(-4122)   # xlPasteFormats

$ws.Cells.Item(118, 1).Value2 = "89.pdf"
$ws.Cells.Item(118, 2).Value2 = "The Ladies and Gentlemen's Complete Letter-Writer Containing Familiar Letters in the Most Common Occasions in Life. Also, a Variety of Elegant Letters for the Direction and Embellishment of Style, on Business, Duty, Amusement, Love, Courtship, Marriage, Friendship, and Other Subjects with Directions for Writing Letters, and the Proper Forms of Address"
$ws.Cells.Item(118, 3).Value2 = "Anon."
$ws.Cells.Item(118, 4).Value2 = 1797

# --- New "Pays" column (F) for every record, old and new -----------------
$ws.Range("F2:F118").Value2 = "British Isles"

# --- Match the author's final selection in the sheet ----------------------
$ws.Range("I116").Select() | Out-Null
